$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix A27: was stored as text "71277628", should be numeric 71277628
$ws.Range("A27").Value = 71277628

# 2. Append new row 28 for payment 71277628 (Cash) 2025-08-18T16:54:54
#    A28 keeps the "raw" text form (mirrors how A27 originally looked before
#    being normalised to a number), so force text via a leading apostrophe,
#    then drop back to the default "Normal" style so no stray quote-prefix
#    formatting is left behind on the cell.
$ws.Range("A28").Value = "'71277628"
$ws.Range("A28").Style = "Normal"
$ws.Range("C28").Value = "Cash"
$ws.Range("D28").Value = "2025-08-18T16:54:54"
$ws.Range("E28").Value = 76
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 76
